$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diciembre")

$ws.Range("R4").Value = "⬇ + 1d 14h"
$ws.Range("C9").Value = "🟢 11h 17m"
$ws.Range("B12").Value = "🔵 5d 6h"
$ws.Range("E12").Value = "🟡 21h 41m"
$ws.Range("R12").Value = "⬇ + 5h 40m"
$ws.Range("C14").Value = "🟢 8h 1m"
$ws.Range("O14").Value = "⬇ + 7d 19h"
$ws.Range("P14").Value = "⬆ - 23h 38m"
